$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-F (Id, Date, Time, League, Home, Away) hold text values, including
# date-/time-shaped strings like "07/11/2024". Force text format before writing so
# Excel does not auto-convert them into date serials, then restore the default
# style so the cells end up unstyled (matching the source data, which carries no
# explicit cell style).
$ws.Range("A2:F5").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "boe4DyiR"
$ws.Range("B2").Value = "07/11/2024"
$ws.Range("C2").Value = "12:00"
$ws.Range("D2").Value = "EGYPT - PREMIER LEAGUE"
$ws.Range("E2").Value = "Al Ittihad"
$ws.Range("F2").Value = "Haras El Hodood"
$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 3.15
$ws.Range("I2").Value = 5.1
$ws.Range("J2").Value = 2.32
$ws.Range("K2").Value = 2.05
$ws.Range("L2").Value = 5.3
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 6.1
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.72
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 1.45
$ws.Range("T2").Value = 2.57
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.72
$ws.Range("W2").Value = 5.6
$ws.Range("X2").Value = 7.4
$ws.Range("Y2").Value = 8.25
$ws.Range("Z2").Value = 14
$ws.Range("AA2").Value = 15.5
$ws.Range("AB2").Value = 32
$ws.Range("AC2").Value = 6.1
$ws.Range("AD2").Value = 6.2
$ws.Range("AE2").Value = 17
$ws.Range("AF2").Value = 100
$ws.Range("AG2").Value = 900
$ws.Range("AH2").Value = 11.5
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 16.5
$ws.Range("AK2").Value = 100
$ws.Range("AL2").Value = 60
$ws.Range("AM2").Value = 60
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 8.75
$ws.Range("AP2").Value = 18.5
$ws.Range("AQ2").Value = 32
$ws.Range("AR2").Value = 65
$ws.Range("AS2").Value = 250
$ws.Range("AT2").Value = 2.57
$ws.Range("AU2").Value = 7.4
$ws.Range("AV2").Value = 70
$ws.Range("AW2").Value = 6.7
$ws.Range("AX2").Value = 30
$ws.Range("AY2").Value = 35
$ws.Range("AZ2").Value = 200
$ws.Range("BA2").Value = 250
$ws.Range("BB2").Value = 450
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51

# Row 3
$ws.Range("A3").Value = "jLidFcME"
$ws.Range("B3").Value = "07/11/2024"
$ws.Range("C3").Value = "12:00"
$ws.Range("D3").Value = "EGYPT - PREMIER LEAGUE"
$ws.Range("E3").Value = "Ghazl El Mahallah"
$ws.Range("F3").Value = "El Gaish"
$ws.Range("G3").Value = 3.05
$ws.Range("H3").Value = 2.62
$ws.Range("I3").Value = 2.72
$ws.Range("J3").Value = 3.7
$ws.Range("K3").Value = 1.85
$ws.Range("L3").Value = 3.35
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 5.2
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.42
$ws.Range("Q3").Value = 2.42
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 1.55
$ws.Range("T3").Value = 2.35
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 7.3
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 10.75
$ws.Range("Z3").Value = 45
$ws.Range("AA3").Value = 32
$ws.Range("AB3").Value = 45
$ws.Range("AC3").Value = 5.2
$ws.Range("AD3").Value = 5.2
$ws.Range("AE3").Value = 14.5
$ws.Range("AF3").Value = 80
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 6.8
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 10
$ws.Range("AK3").Value = 35
$ws.Range("AL3").Value = 27
$ws.Range("AM3").Value = 40
$ws.Range("AN3").Value = 4.85
$ws.Range("AO3").Value = 18
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 90
$ws.Range("AR3").Value = 150
$ws.Range("AS3").Value = 350
$ws.Range("AT3").Value = 2.32
$ws.Range("AU3").Value = 6.8
$ws.Range("AV3").Value = 65
$ws.Range("AW3").Value = 4.5
$ws.Range("AX3").Value = 15.5
$ws.Range("AY3").Value = 24
$ws.Range("AZ3").Value = 75
$ws.Range("BA3").Value = 120
$ws.Range("BB3").Value = 350
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51

# Row 4
$ws.Range("A4").Value = "YyDJubM9"
$ws.Range("B4").Value = "07/11/2024"
$ws.Range("C4").Value = "11:30"
$ws.Range("D4").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E4").Value = "Al Qadisiya"
$ws.Range("F4").Value = "Al Feiha"
$ws.Range("G4").Value = 1.4
$ws.Range("H4").Value = 4.5
$ws.Range("I4").Value = 6.25
$ws.Range("J4").Value = 1.91
$ws.Range("K4").Value = 2.3
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 1.02
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 6.5
$ws.Range("X4").Value = 6.5
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 9
$ws.Range("AA4").Value = 13
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 9
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 1500
$ws.Range("AH4").Value = 15
$ws.Range("AI4").Value = 34
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 81
$ws.Range("AL4").Value = 51
$ws.Range("AM4").Value = 51
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 7
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 19
$ws.Range("AR4").Value = 41
$ws.Range("AS4").Value = 151
$ws.Range("AT4").Value = 3
$ws.Range("AU4").Value = 9.5
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 8.5
$ws.Range("AX4").Value = 41
$ws.Range("AY4").Value = 41
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 151
$ws.Range("BB4").Value = 500
$ws.Range("BC4").Value = 81
$ws.Range("BD4").Value = 81

# Row 5
$ws.Range("A5").Value = "8fR1hy6F"
$ws.Range("B5").Value = "07/11/2024"
$ws.Range("C5").Value = "12:00"
$ws.Range("D5").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E5").Value = "Al Kholood"
$ws.Range("F5").Value = "Al Shabab"
$ws.Range("G5").Value = 3.75
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 1.9
$ws.Range("J5").Value = 4.33
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 2.5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.05
$ws.Range("R5").Value = 1.75
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.83
$ws.Range("W5").Value = 11
$ws.Range("X5").Value = 19
$ws.Range("Y5").Value = 15
$ws.Range("Z5").Value = 41
$ws.Range("AA5").Value = 34
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 6.5
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 51
$ws.Range("AG5").Value = 700
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = 8.5
$ws.Range("AJ5").Value = 9
$ws.Range("AK5").Value = 15
$ws.Range("AL5").Value = 17
$ws.Range("AM5").Value = 29
$ws.Range("AN5").Value = 6
$ws.Range("AO5").Value = 21
$ws.Range("AP5").Value = 29
$ws.Range("AQ5").Value = 67
$ws.Range("AR5").Value = 101
$ws.Range("AS5").Value = 300
$ws.Range("AT5").Value = 2.63
$ws.Range("AU5").Value = 8.5
$ws.Range("AV5").Value = 51
$ws.Range("AW5").Value = 4
$ws.Range("AX5").Value = 10
$ws.Range("AY5").Value = 21
$ws.Range("AZ5").Value = 34
$ws.Range("BA5").Value = 51
$ws.Range("BB5").Value = 151
$ws.Range("BC5").Value = 81
$ws.Range("BD5").Value = 81

# Restore default (unstyled) formatting on the text columns now that the values
# are safely stored as text, so cells end up with no explicit style, as in the source.
$ws.Range("A2:F5").Style = "Normal"
